# Rename the sheet and update the active cell selection to match the
# workbook as re-uploaded by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Tabelle1" -> "P4TSSOP8_505_32X6_1M_85HS120_F1"
$ws.Name = "P4TSSOP8_505_32X6_1M_85HS120_F1"

# Move the selected/active cell from Q11 to G16
$ws.Range("G16").Select()
